$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 420
$ws.Cells.Item(420, 1).Value = 419
$ws.Cells.Item(420, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(420, 3).Value = "1:25 AM"
$ws.Cells.Item(420, 4).Value = "W61612"
$ws.Cells.Item(420, 5).Value = "Liverpool"
$ws.Cells.Item(420, 6).Value = "(LPL)"
$ws.Cells.Item(420, 7).Value = "Wizz Air "
$ws.Cells.Item(420, 8).Value = "A321"
$ws.Cells.Item(420, 9).Value = "(HA-LXL)"
$ws.Cells.Item(420, 10).Value = "1:57 AM"
$ws.Cells.Item(420, 12).Value = "0 hours, 32 minutes"

# Row 421
$ws.Cells.Item(421, 1).Value = 420
$ws.Cells.Item(421, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(421, 3).Value = "2:25 AM"
$ws.Cells.Item(421, 4).Value = "W61774"
$ws.Cells.Item(421, 5).Value = "Reykjavik"
$ws.Cells.Item(421, 6).Value = "(KEF)"
$ws.Cells.Item(421, 7).Value = "Wizz Air "
$ws.Cells.Item(421, 8).Value = "A320"
$ws.Cells.Item(421, 9).Value = "(HA-LYH)"
$ws.Cells.Item(421, 10).Value = "2:01 AM"
$ws.Cells.Item(421, 12).Value = "0 hours, -24 minutes"

# Row 422
$ws.Cells.Item(422, 1).Value = 421
$ws.Cells.Item(422, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(422, 3).Value = "9:25 AM"
$ws.Cells.Item(422, 4).Value = "FR6102"
$ws.Cells.Item(422, 5).Value = "Stockholm"
$ws.Cells.Item(422, 6).Value = "(ARN)"
$ws.Cells.Item(422, 7).Value = "Ryanair "
$ws.Cells.Item(422, 8).Value = "B738"
$ws.Cells.Item(422, 9).Value = "(SP-RKM)"
$ws.Cells.Item(422, 10).Value = "9:11 AM"
$ws.Cells.Item(422, 12).Value = "0 hours, -14 minutes"

# Row 423
$ws.Cells.Item(423, 1).Value = 422
$ws.Cells.Item(423, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(423, 3).Value = "9:30 AM"
$ws.Cells.Item(423, 4).Value = "DY1030"
$ws.Cells.Item(423, 5).Value = "Bergen"
$ws.Cells.Item(423, 6).Value = "(BGO)"
$ws.Cells.Item(423, 7).Value = "Norwegian "
$ws.Cells.Item(423, 8).Value = "B738"
$ws.Cells.Item(423, 9).Value = "(LN-DYM)"
$ws.Cells.Item(423, 10).Value = "9:27 AM"
$ws.Cells.Item(423, 12).Value = "0 hours, -3 minutes"

# Row 424
$ws.Cells.Item(424, 1).Value = 423
$ws.Cells.Item(424, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(424, 3).Value = "9:45 AM"
$ws.Cells.Item(424, 4).Value = "W61744"
$ws.Cells.Item(424, 5).Value = "Oslo"
$ws.Cells.Item(424, 6).Value = "(TRF)"
$ws.Cells.Item(424, 7).Value = "Wizz Air "
$ws.Cells.Item(424, 8).Value = "A320"
$ws.Cells.Item(424, 9).Value = "(HA-LYO)"
$ws.Cells.Item(424, 10).Value = "9:45 AM"
$ws.Cells.Item(424, 12).Value = "0 hours, 0 minutes"

# Row 425
$ws.Cells.Item(425, 1).Value = 424
$ws.Cells.Item(425, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(425, 3).Value = "10:10 AM"
$ws.Cells.Item(425, 4).Value = "FR3687"
$ws.Cells.Item(425, 5).Value = "Billund"
$ws.Cells.Item(425, 6).Value = "(BLL)"
$ws.Cells.Item(425, 7).Value = "Ryanair "
$ws.Cells.Item(425, 8).Value = "B738"
$ws.Cells.Item(425, 9).Value = "(SP-RSW)"
$ws.Cells.Item(425, 10).Value = "9:50 AM"
$ws.Cells.Item(425, 12).Value = "0 hours, -20 minutes"

# Row 426
$ws.Cells.Item(426, 1).Value = 425
$ws.Cells.Item(426, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(426, 3).Value = "10:35 AM"
$ws.Cells.Item(426, 4).Value = "W61642"
$ws.Cells.Item(426, 5).Value = "Eindhoven"
$ws.Cells.Item(426, 6).Value = "(EIN)"
$ws.Cells.Item(426, 7).Value = "Wizz Air "
$ws.Cells.Item(426, 8).Value = "A21N"
$ws.Cells.Item(426, 9).Value = "(9H-WBU)"
$ws.Cells.Item(426, 10).Value = "10:23 AM"
$ws.Cells.Item(426, 12).Value = "0 hours, -12 minutes"

# Row 427
$ws.Cells.Item(427, 1).Value = 426
$ws.Cells.Item(427, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(427, 3).Value = "11:30 AM"
$ws.Cells.Item(427, 4).Value = "FR2374"
$ws.Cells.Item(427, 5).Value = "London"
$ws.Cells.Item(427, 6).Value = "(STN)"
$ws.Cells.Item(427, 7).Value = "Ryanair "
$ws.Cells.Item(427, 8).Value = "B738"
$ws.Cells.Item(427, 9).Value = "(SP-RKQ)"
$ws.Cells.Item(427, 10).Value = "12:07 PM"
$ws.Cells.Item(427, 12).Value = "0 hours, 37 minutes"

# Row 428
$ws.Cells.Item(428, 1).Value = 427
$ws.Cells.Item(428, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(428, 3).Value = "11:30 AM"
$ws.Cells.Item(428, 4).Value = "W61602"
$ws.Cells.Item(428, 5).Value = "London"
$ws.Cells.Item(428, 6).Value = "(LTN)"
$ws.Cells.Item(428, 7).Value = "Wizz Air "
$ws.Cells.Item(428, 8).Value = "A320"
$ws.Cells.Item(428, 9).Value = "(HA-LYH)"
$ws.Cells.Item(428, 10).Value = "11:21 AM"
$ws.Cells.Item(428, 12).Value = "0 hours, -9 minutes"

# Row 429
$ws.Cells.Item(429, 1).Value = 428
$ws.Cells.Item(429, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(429, 3).Value = "12:05 PM"
$ws.Cells.Item(429, 4).Value = "FR3286"
$ws.Cells.Item(429, 5).Value = "Leeds"
$ws.Cells.Item(429, 6).Value = "(LBA)"
$ws.Cells.Item(429, 7).Value = "Ryanair "
$ws.Cells.Item(429, 8).Value = "B738"
$ws.Cells.Item(429, 9).Value = "(SP-RSO)"
$ws.Cells.Item(429, 10).Value = "11:43 AM"
$ws.Cells.Item(429, 12).Value = "0 hours, -22 minutes"

# Row 430
$ws.Cells.Item(430, 1).Value = 429
$ws.Cells.Item(430, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(430, 3).Value = "12:35 PM"
$ws.Cells.Item(430, 4).Value = "FR6845"
$ws.Cells.Item(430, 5).Value = "Copenhagen"
$ws.Cells.Item(430, 6).Value = "(CPH)"
$ws.Cells.Item(430, 7).Value = "Ryanair "
$ws.Cells.Item(430, 8).Value = "B738"
$ws.Cells.Item(430, 9).Value = "(SP-RKM)"
$ws.Cells.Item(430, 10).Value = "12:15 PM"
$ws.Cells.Item(430, 12).Value = "0 hours, -20 minutes"

# Row 431
$ws.Cells.Item(431, 1).Value = 430
$ws.Cells.Item(431, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(431, 3).Value = "12:35 PM"
$ws.Cells.Item(431, 4).Value = "LH1642"
$ws.Cells.Item(431, 5).Value = "Munich"
$ws.Cells.Item(431, 6).Value = "(MUC)"
$ws.Cells.Item(431, 7).Value = "Lufthansa "
$ws.Cells.Item(431, 8).Value = "CRJ9"
$ws.Cells.Item(431, 9).Value = "(D-ACNI)"
$ws.Cells.Item(431, 10).Value = "12:21 PM"
$ws.Cells.Item(431, 12).Value = "0 hours, -14 minutes"

# Row 432
$ws.Cells.Item(432, 1).Value = 431
$ws.Cells.Item(432, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(432, 3).Value = "12:45 PM"
$ws.Cells.Item(432, 4).Value = "FR6126"
$ws.Cells.Item(432, 5).Value = "Malta"
$ws.Cells.Item(432, 6).Value = "(MLA)"
$ws.Cells.Item(432, 7).Value = "Ryanair "
$ws.Cells.Item(432, 8).Value = "B738"
$ws.Cells.Item(432, 9).Value = "(SP-RSL)"
$ws.Cells.Item(432, 10).Value = "1:04 PM"
$ws.Cells.Item(432, 12).Value = "0 hours, 19 minutes"

# Row 433
$ws.Cells.Item(433, 1).Value = 432
$ws.Cells.Item(433, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(433, 3).Value = "1:35 PM"
$ws.Cells.Item(433, 4).Value = "LH1376"
$ws.Cells.Item(433, 5).Value = "Frankfurt"
$ws.Cells.Item(433, 6).Value = "(FRA)"
$ws.Cells.Item(433, 7).Value = "Lufthansa "
$ws.Cells.Item(433, 8).Value = "CRJ9"
$ws.Cells.Item(433, 9).Value = "(D-ACKI)"
$ws.Cells.Item(433, 10).Value = "1:45 PM"
$ws.Cells.Item(433, 12).Value = "0 hours, 10 minutes"

# Row 434
$ws.Cells.Item(434, 1).Value = 433
$ws.Cells.Item(434, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(434, 3).Value = "2:10 PM"
$ws.Cells.Item(434, 4).Value = "LO3837"
$ws.Cells.Item(434, 5).Value = "Warsaw"
$ws.Cells.Item(434, 6).Value = "(WAW)"
$ws.Cells.Item(434, 7).Value = "LOT "
$ws.Cells.Item(434, 8).Value = "E75S"
$ws.Cells.Item(434, 9).Value = "(SP-LIB)"
$ws.Cells.Item(434, 10).Value = "1:54 PM"
$ws.Cells.Item(434, 12).Value = "0 hours, -16 minutes"

# Row 435
$ws.Cells.Item(435, 1).Value = 434
$ws.Cells.Item(435, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(435, 3).Value = "2:55 PM"
$ws.Cells.Item(435, 4).Value = "DY1054"
$ws.Cells.Item(435, 5).Value = "Oslo"
$ws.Cells.Item(435, 6).Value = "(OSL)"
$ws.Cells.Item(435, 7).Value = "Norwegian "
$ws.Cells.Item(435, 8).Value = "B738"
$ws.Cells.Item(435, 9).Value = "(LN-DYU)"
$ws.Cells.Item(435, 10).Value = "3:07 PM"
$ws.Cells.Item(435, 12).Value = "0 hours, 12 minutes"

# Row 436
$ws.Cells.Item(436, 1).Value = 435
$ws.Cells.Item(436, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(436, 3).Value = "2:55 PM"
$ws.Cells.Item(436, 4).Value = "W61740"
$ws.Cells.Item(436, 5).Value = "Malmo"
$ws.Cells.Item(436, 6).Value = "(MMX)"
$ws.Cells.Item(436, 7).Value = "Wizz Air "
$ws.Cells.Item(436, 8).Value = "A321"
$ws.Cells.Item(436, 9).Value = "(HA-LTB)"
$ws.Cells.Item(436, 10).Value = "2:36 PM"
$ws.Cells.Item(436, 12).Value = "0 hours, -19 minutes"

# Row 437
$ws.Cells.Item(437, 1).Value = 436
$ws.Cells.Item(437, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(437, 3).Value = "3:05 PM"
$ws.Cells.Item(437, 4).Value = "SK759"
$ws.Cells.Item(437, 5).Value = "Copenhagen"
$ws.Cells.Item(437, 6).Value = "(CPH)"
$ws.Cells.Item(437, 7).Value = "SAS "
$ws.Cells.Item(437, 8).Value = "A20N"
$ws.Cells.Item(437, 9).Value = "(EI-SIA)"
$ws.Cells.Item(437, 10).Value = "2:54 PM"
$ws.Cells.Item(437, 12).Value = "0 hours, -11 minutes"

# Row 438
$ws.Cells.Item(438, 1).Value = 437
$ws.Cells.Item(438, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(438, 3).Value = "3:55 PM"
$ws.Cells.Item(438, 4).Value = "FR8509"
$ws.Cells.Item(438, 5).Value = "Oslo"
$ws.Cells.Item(438, 6).Value = "(TRF)"
$ws.Cells.Item(438, 7).Value = "Ryanair "
$ws.Cells.Item(438, 8).Value = "B738"
$ws.Cells.Item(438, 9).Value = "(SP-RSO)"
$ws.Cells.Item(438, 10).Value = "3:46 PM"
$ws.Cells.Item(438, 12).Value = "0 hours, -9 minutes"

# Row 439
$ws.Cells.Item(439, 1).Value = 438
$ws.Cells.Item(439, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(439, 3).Value = "4:25 PM"
$ws.Cells.Item(439, 4).Value = "FR6110"
$ws.Cells.Item(439, 5).Value = "Manchester"
$ws.Cells.Item(439, 6).Value = "(MAN)"
$ws.Cells.Item(439, 7).Value = "Ryanair "
$ws.Cells.Item(439, 8).Value = "B738"
$ws.Cells.Item(439, 9).Value = "(SP-RSW)"
$ws.Cells.Item(439, 10).Value = "3:58 PM"
$ws.Cells.Item(439, 12).Value = "0 hours, -27 minutes"
